$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as plain text so that
# values such as "314.57" or "6.62" are not reinterpreted as numbers
# by Excel's automatic type coercion on Range.Value assignment.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

# Mapping of row number -> (new Price text or $null if unchanged, new Volume(1h) text)
$updates = @{
    2  = @{ D = "40.805.48"; E = "  -2.03%  " }
    3  = @{ D = "2.399.95";  E = "  -2.95%  " }
    4  = @{ D = $null;       E = "  +0.36%  " }
    5  = @{ D = "314.57";    E = "  -0.92%  " }
    6  = @{ D = "87.69";     E = "  -5.44%  " }
    7  = @{ D = "0.535";     E = "  -3.06%  " }
    8  = @{ D = $null;       E = "  +0.08%  " }
    9  = @{ D = "0.490";     E = "  -5.09%  " }
    10 = @{ D = $null;       E = "  -2.88%  " }
    11 = @{ D = "31.36";     E = "  -5.16%  " }
    12 = @{ D = $null;       E = "  -1.83%  " }
    13 = @{ D = "2.780.18";  E = "  -2.68%  " }
    14 = @{ D = "6.62";      E = "  -4.18%  " }
    15 = @{ D = "15.47";     E = "  -1.84%  " }
    16 = @{ D = "2.418.30";  E = "  -2.93%  " }
    17 = @{ D = "0.764";     E = "  -3.41%  " }
    18 = @{ D = "40.723.51"; E = "  -2.14%  " }
    19 = @{ D = "0.0₃0912";  E = "  -3.80%  " }
    20 = @{ D = "6.19";      E = "  -4.16%  " }
    21 = @{ D = "71.17";     E = "  -0.21%  " }
    22 = @{ D = "10.77";     E = "  -4.53%  " }
    23 = @{ D = "231.20";    E = "  -3.43%  " }
    24 = @{ D = $null;       E = "  -2.84%  " }
    25 = @{ D = $null;       E = "  +0.21%  " }
    26 = @{ D = "1.83";      E = "  -5.02%  " }
    27 = @{ D = "23.75";     E = "  -3.82%  " }
    28 = @{ D = $null;       E = "  -2.67%  " }
    29 = @{ D = "9.44";      E = "  -3.97%  " }
    30 = @{ D = "33.89";     E = "  -6.42%  " }
    31 = @{ D = "156.84";    E = "  -1.97%  " }
    32 = @{ D = $null;       E = "  +0.03%  " }
    33 = @{ D = "5.19";      E = "  -5.84%  " }
    34 = @{ D = "0.0731";    E = "  -4.86%  " }
    35 = @{ D = $null;       E = "  -5.55%  " }
    36 = @{ D = "2.86";      E = "  -2.21%  " }
    37 = @{ D = "16.21";     E = "  -6.48%  " }
    38 = @{ D = $null;       E = "  -2.21%  " }
    39 = @{ D = "1.75";      E = "  -5.57%  " }
    40 = @{ D = "0.0992";    E = "  -4.47%  " }
    41 = @{ D = "3.82";      E = "  -4.31%  " }
    42 = @{ D = $null;       E = "  -6.98%  " }
    43 = @{ D = "1.975.70";  E = "  -0.47%  " }
    44 = @{ D = "18.29";     E = "  -2.28%  " }
    45 = @{ D = "0.0271";    E = "  -5.18%  " }
    46 = @{ D = "2.82";      E = "  -5.69%  " }
    47 = @{ D = "9.28";      E = "  -0.27%  " }
    48 = @{ D = "2.647.66";  E = "  -2.44%  " }
    49 = @{ D = "93.35";     E = "  -4.06%  " }
    50 = @{ D = "72.64";     E = "  -1.90%  " }
    51 = @{ D = "50.80";     E = "  -2.70%  " }
}

foreach ($row in $updates.Keys) {
    $entry = $updates[$row]
    if ($null -ne $entry.D) {
        $ws.Cells.Item($row, 4).Value = $entry.D
    }
    $ws.Cells.Item($row, 5).Value = $entry.E
}

# Restore the original (default) cell style now that the text values are
# safely stored, so no stray number-format style is left behind on cells
# that previously had no explicit style.
$textRange.Style = "Normal"
